$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (NCTId), shifting existing columns C:L to D:M.
$ws.Columns("C:C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "statut_name"

# Populate the new "statut_name" column based on each row's statut_label (column B).
$map = @{
    "noir"   = "pas de résultat ni de publication"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    if ($map.ContainsKey($label)) {
        $ws.Cells.Item($r, 3).Value = $map[$label]
    }
}
